$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Metadata sheet: Version, Date, Count updates. Date/Count must
#    stay plain text (not get auto-converted to a date serial / a
#    number), so we briefly force a text number-format while writing
#    them, then restore the original plain body style by pasting the
#    format from an already-correctly-styled text cell on the sheet.
# ------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")

$wsMeta.Range("B3").Value = "1.0.2"

$wsMeta.Range("B8").NumberFormat = "@"
$wsMeta.Range("B8").Value = "2025-09-22"
$wsMeta.Range("A7").Copy()
$wsMeta.Range("B8").PasteSpecial(-4122)

$wsMeta.Range("B22").NumberFormat = "@"
$wsMeta.Range("B22").Value = "18"
$wsMeta.Range("A21").Copy()
$wsMeta.Range("B22").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 2) Rename the old "Concepts" sheet to "Properties", then add a
#    fresh "Concepts" sheet right after it that will carry the actual
#    concept rows (the 16 original concepts plus 2 new ones).
# ------------------------------------------------------------------
$wsOldConcepts = $wb.Worksheets.Item("Concepts")
$wsOldConcepts.Name = "Properties"

$wsConcepts = $wb.Worksheets.Add($null, $wsOldConcepts)
$wsConcepts.Name = "Concepts"

# Clone the header / body formatting from the (still concept-shaped)
# Properties sheet onto the new Concepts sheet before the Properties
# sheet's own rows get rewritten below.
$wsOldConcepts.Range("A1:D1").Copy()
$wsConcepts.Range("A1:D1").PasteSpecial(-4122)

$wsOldConcepts.Range("A2:D2").Copy()
$wsConcepts.Range("A2:D19").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 3) Rewrite the Properties sheet: it now documents the two FHIR
#    concept properties (status, effectiveDate) instead of concepts.
# ------------------------------------------------------------------
$wsOldConcepts.Rows("4:17").Delete()

$wsOldConcepts.Range("A1").Value = "Code"
$wsOldConcepts.Range("B1").Value = "Uri"
$wsOldConcepts.Range("C1").Value = "Description"
$wsOldConcepts.Range("D1").Value = "Type"

$wsOldConcepts.Range("A2").Value = "status"
$wsOldConcepts.Range("B2").Value = "http://hl7.org/fhir/concept-properties#status"
$wsOldConcepts.Range("C2").Value = "A property that indicates the status of the concept. One of active, experimental, deprecated, or retired."
$wsOldConcepts.Range("D2").Value = "code"

$wsOldConcepts.Range("A3").Value = "effectiveDate"
$wsOldConcepts.Range("B3").Value = "http://hl7.org/fhir/concept-properties#effectiveDate"
$wsOldConcepts.Range("C3").Value = "The date at which the concept status was last changed."
$wsOldConcepts.Range("D3").Value = "dateTime"

# ------------------------------------------------------------------
# 4) Populate the new Concepts sheet with the full concept list: the
#    original 16 rows plus 2 new ones ("acute-ambulant" and
#    "extended-care-responsibility"), inserted before "other".
# ------------------------------------------------------------------
$wsConcepts.Range("A1").Value = "Level"
$wsConcepts.Range("B1").Value = "Code"
$wsConcepts.Range("C1").Value = "Display"
$wsConcepts.Range("D1").Value = "Definition"

$concepts = @(
  @("alcohol-and-drug-treatment", "Alcohol and drug treatment"),
  @("assistive-devices", "Assistive technology"),
  @("carecoordination", "Care Coordination"),
  @("decease", "Decease"),
  @("discharge", "Discharge"),
  @("examination-results", "Examination Results"),
  @("healthcare", "Healthcare"),
  @("home-care-assessment", "Home care assessment"),
  @("medicine", "Medicine"),
  @("nursing", "Nursing"),
  @("outpatient", "Outpatient"),
  @("psychiatry-social-disability", "Psychiatry, Social, Disability"),
  @("regarding-referral", "Regarding Referral"),
  @("telemedicine", "Telemedicine"),
  @("training", "Training"),
  @("acute-ambulant", "Acute ambulant"),
  @("extended-care-responsibility", "Extended care responsibility"),
  @("other", "Other")
)

$row = 2
foreach ($c in $concepts) {
    $wsConcepts.Range("B$row").Value = $c[0]
    $wsConcepts.Range("C$row").Value = $c[1]
    $row = $row + 1
}

# The "Level" column is always literal text "1" (not a number). Force
# text storage the same way as above, then restore the plain body
# style by pasting the (already correct) format from column B.
$wsConcepts.Range("A2:A19").NumberFormat = "@"
$wsConcepts.Range("A2:A19").Value = "1"
$wsConcepts.Range("B2:B19").Copy()
$wsConcepts.Range("A2:A19").PasteSpecial(-4122)
